$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H46").Value = 1878.1818
$ws.Range("I46").Value = 980
$ws.Range("J46").Value = 1968
$ws.Range("K46").Value = 2940
$ws.Range("L46").Value = 5904
$ws.Range("M46").Value = -2821
$ws.Range("N46").Value = -6142
$ws.Range("H60").Value = 1878.1818
$ws.Range("I60").Value = 980
$ws.Range("J60").Value = 1968
$ws.Range("K60").Value = 2940
$ws.Range("L60").Value = 5904
$ws.Range("M60").Value = -2456
$ws.Range("N60").Value = -6872
$ws.Range("H87").Value = 22425.455
$ws.Range("J87").Value = 22425.455
$ws.Range("L87").Value = 22425.455
$ws.Range("N87").Value = -24921.455
$ws.Range("H90").Value = 22425.455
$ws.Range("J90").Value = 22425.455
$ws.Range("L90").Value = 67276.36500000001
$ws.Range("N90").Value = -79756.36500000001
$ws.Range("H138").Value = 2610.7903
$ws.Range("I138").Value = 1187.85
$ws.Range("J138").Value = 5197.9546
$ws.Range("K138").Value = 3563.55
$ws.Range("L138").Value = 15593.8638
$ws.Range("M138").Value = 1576.45
$ws.Range("N138").Value = -25873.8638

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3468.899
$ws.Range("I32").Value = 2424.0342
$ws.Range("K32").Value = 2424.0342
$ws.Range("M32").Value = -2137.0342
$ws.Range("H61").Value = 1056.1578
$ws.Range("I61").Value = 1056.1578
$ws.Range("K61").Value = 1056.1578
$ws.Range("M61").Value = -844.1578
$ws.Range("H74").Value = 1790
$ws.Range("I74").Value = 547.1
$ws.Range("J74").Value = 4275.8
$ws.Range("K74").Value = 547.1
$ws.Range("L74").Value = 4275.8
$ws.Range("M74").Value = 326.9
$ws.Range("N74").Value = -6023.8
$ws.Range("H77").Value = 1790
$ws.Range("I77").Value = 547.1
$ws.Range("J77").Value = 4275.8
$ws.Range("K77").Value = 2735.5
$ws.Range("L77").Value = 21379
$ws.Range("M77").Value = 1632.5
$ws.Range("N77").Value = -30115
$ws.Range("H122").Value = 5265404.5
$ws.Range("I122").Value = 6668078
$ws.Range("J122").Value = 5378.5
$ws.Range("K122").Value = 20004234
$ws.Range("L122").Value = 16135.5
$ws.Range("M122").Value = -20001784
$ws.Range("N122").Value = -21035.5
$ws.Range("H136").Value = 1056.1578
$ws.Range("I136").Value = 1056.1578
$ws.Range("K136").Value = 3168.4734
$ws.Range("M136").Value = -618.4733999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 295.42856
$ws.Range("I22").Value = 281.33334
$ws.Range("J22").Value = 380
$ws.Range("K22").Value = 281.33334
$ws.Range("L22").Value = 380
$ws.Range("M22").Value = -108.33334
$ws.Range("N22").Value = -726
$ws.Range("H105").Value = 2734.9546
$ws.Range("I105").Value = 3197.6924
$ws.Range("J105").Value = 2066.5557
$ws.Range("K105").Value = 3197.6924
$ws.Range("L105").Value = 2066.5557
$ws.Range("M105").Value = -1450.6924
$ws.Range("N105").Value = -5560.5557

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1155.3962
$ws.Range("I58").Value = 1009.9048
$ws.Range("J58").Value = 1710.909
$ws.Range("K58").Value = 1009.9048
$ws.Range("L58").Value = 1710.909
$ws.Range("M58").Value = -806.9048
$ws.Range("N58").Value = -2116.909
$ws.Range("H136").Value = 1155.3962
$ws.Range("I136").Value = 1009.9048
$ws.Range("J136").Value = 1710.909
$ws.Range("K136").Value = 3029.7144
$ws.Range("L136").Value = 5132.727000000001
$ws.Range("M136").Value = -479.7143999999998
$ws.Range("N136").Value = -10232.727

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 678.2222
$ws.Range("I33").Value = 100
$ws.Range("J33").Value = 750.5
$ws.Range("K33").Value = 600
$ws.Range("L33").Value = 4503
$ws.Range("M33").Value = -317
$ws.Range("N33").Value = -5069
$ws.Range("H34").Value = 2218.2222
$ws.Range("I34").Value = 1546
$ws.Range("J34").Value = 2756
$ws.Range("K34").Value = 4638
$ws.Range("L34").Value = 8268
$ws.Range("M34").Value = -4554
$ws.Range("N34").Value = -8436
$ws.Range("H39").Value = 1912.8
$ws.Range("J39").Value = 1912.8
$ws.Range("L39").Value = 5738.4
$ws.Range("N39").Value = -6326.4
$ws.Range("H49").Value = 3000
$ws.Range("J49").Value = 3000
$ws.Range("L49").Value = 9000
$ws.Range("N49").Value = -9312
$ws.Range("H55").Value = 2113.5
$ws.Range("I55").Value = 973.7143
$ws.Range("K55").Value = 2921.1429
$ws.Range("M55").Value = -2744.1429
$ws.Range("H105").Value = 7978.5713
$ws.Range("J105").Value = 7978.5713
$ws.Range("L105").Value = 23935.7139
$ws.Range("N105").Value = -29177.7139

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5357.9165
$ws.Range("I70").Value = 5185.7144
$ws.Range("J70").Value = 5599
$ws.Range("K70").Value = 5185.7144
$ws.Range("L70").Value = 5599
$ws.Range("M70").Value = -4915.7144
$ws.Range("N70").Value = -6139
$ws.Range("H73").Value = 5357.9165
$ws.Range("I73").Value = 5185.7144
$ws.Range("J73").Value = 5599
$ws.Range("K73").Value = 5185.7144
$ws.Range("L73").Value = 5599
$ws.Range("M73").Value = -4249.7144
$ws.Range("N73").Value = -7471
$ws.Range("H80").Value = 2514.2856
$ws.Range("I80").Value = 2500
$ws.Range("J80").Value = 2533.3333
$ws.Range("K80").Value = 2500
$ws.Range("L80").Value = 2533.3333
$ws.Range("M80").Value = -1502
$ws.Range("N80").Value = -4529.3333
$ws.Range("H83").Value = 2514.2856
$ws.Range("I83").Value = 2500
$ws.Range("J83").Value = 2533.3333
$ws.Range("K83").Value = 12500
$ws.Range("L83").Value = 12666.6665
$ws.Range("M83").Value = -7508
$ws.Range("N83").Value = -22650.6665

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 1334.3334
$ws.Range("I55").Value = 1001
$ws.Range("J55").Value = 1501
$ws.Range("K55").Value = 1001
$ws.Range("L55").Value = 1501
$ws.Range("M55").Value = -828
$ws.Range("N55").Value = -1847
$ws.Range("H127").Value = 38514.547
$ws.Range("J127").Value = 38514.547
$ws.Range("L127").Value = 38514.547
$ws.Range("N127").Value = -48434.547
$ws.Range("H132").Value = 1695.7446
$ws.Range("I132").Value = 1464.1666
$ws.Range("J132").Value = 3641
$ws.Range("K132").Value = 4392.4998
$ws.Range("L132").Value = 10923
$ws.Range("M132").Value = -1862.4998
$ws.Range("N132").Value = -15983

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 288.80646
$ws.Range("I136").Value = 265.1
$ws.Range("J136").Value = 1000
$ws.Range("K136").Value = 795.3000000000001
$ws.Range("L136").Value = 3000
$ws.Range("M136").Value = 1754.7
$ws.Range("N136").Value = -8100
